$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

$ws.Range("A" + $row + ":D" + $row).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-01-15"
$ws.Cells.Item($row, 2).Value = "08:57:13"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 5).Value = 126741
$ws.Cells.Item($row, 6).Value = 143481
$ws.Cells.Item($row, 7).Value = 169162
$ws.Cells.Item($row, 8).Value = 155426
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142746
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192963
$ws.Cells.Item($row, 14).Value = 115381
$ws.Cells.Item($row, 15).Value = 45960
$ws.Cells.Item($row, 16).Value = 28500
$ws.Cells.Item($row, 17).Value = 65464
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48436
$ws.Cells.Item($row, 20).Value = -1

$ws.Range("A" + $row + ":D" + $row).ClearFormats()
